$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.394.29'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.822.06'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5259'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3854'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08015'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.115'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.399'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.004'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.422'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").Value = '1.820.88'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001102'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06642'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.034'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("D23").Value = '28.437.92'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.245'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").Value = '2.028.88'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.421'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1107'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.080'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.680'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.676'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07347'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2201'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02340'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.135'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.755'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6321'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.181'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.381'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6127'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.85%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.783'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.981'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.210'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06901'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
